$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix E14:E16 - convert from text (inlineStr) to real numbers
$ws.Range("E14").Value = 543237
$ws.Range("E15").Value = 526371
$ws.Range("E16").Value = 532155

# Row 17
$ws.Range("A17").Value = "11/06/2024 04:44:51"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "HDFCAMC"
$ws.Range("D17").Value = "HDFC Asset Management Company Ltd"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "541729"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = -0.7
$ws.Range("G17").Value = 3764.2
$ws.Range("H17").Value = 141666

# Row 18
$ws.Range("A18").Value = "11/06/2024 04:44:51"
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "MAHABANK"
$ws.Range("D18").Value = "Bank Of Maharashtra"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "532525"
$ws.Range("E18").Style = "Normal"
$ws.Range("F18").Value = -0.12
$ws.Range("G18").Value = 66.09999999999999
$ws.Range("H18").Value = 4373452
